$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DLT Status values
$ws.Range("B2").Value = "Approved"
$ws.Range("B3").Value = "Needs Review"

# Row 5: replace the http-client entry with the new-endpoint entry
$ws.Range("A5").Value = "/src/api/new-endpoint.json"
$ws.Range("B5").Value = "New"
$ws.Range("C5").Value = "API Endpoint"

# Row 6: replace the date-formatter entry with the validator entry
$ws.Range("A6").Value = "/src/utils/validator.js"
$ws.Range("B6").Value = "Approved"
$ws.Range("C6").Value = "Input Validator"

# Remove row 7 entirely (bar-chart.css / Approved / Chart.js)
$ws.Rows.Item(7).Delete()
